# Update parameter descriptions on the "mobility" sheet to clarify that a
# 0.56 loading factor is used to account for empty returning wagons and
# locomotives.
$wb = $excel.ActiveWorkbook

$mobility = $wb.Worksheets.Item("mobility")
$mobility.Range("C7").Value = "Towing capacity of a locomotive (ton). - Using 0,56 loading factor to account for empty returning locomotives."
$mobility.Range("C6").Value = "Loading capacity of a wagon (ton). - Using 0,56 loading factor to account for empty returning wagons."

# Restore the selected cell on each sheet as recorded in the workbook view.
$mobility.Activate()
$mobility.Range("A6").Select()

$infrastructure = $wb.Worksheets.Item("infrastructure")
$infrastructure.Activate()
$infrastructure.Range("C3").Select()

$mobility.Activate()
